$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the "date" column (F) values forward by 15 days for rows 2-7,
# keeping the existing date number format / style intact.
$ws.Range("F2").Value = 44974
$ws.Range("F3").Value = 44973
$ws.Range("F4").Value = 44972
$ws.Range("F5").Value = 44971
$ws.Range("F6").Value = 44970
$ws.Range("F7").Value = 44969
